$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text so numeric-looking strings
# (e.g. "6.00", "245.12") are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "43.071.83"
$ws.Range("E2").Value = "  +5.65%  "
$ws.Range("D3").Value = "2.241.07"
$ws.Range("E3").Value = "  +5.30%  "
$ws.Range("E4").Value = "  -0.37%  "
$ws.Range("D5").Value = "245.12"
$ws.Range("E5").Value = "  +5.49%  "
$ws.Range("E6").Value = "  +4.16%  "
$ws.Range("D7").Value = "75.25"
$ws.Range("E7").Value = "  +11.34%  "
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("E9").Value = "  +9.34%  "
$ws.Range("D10").Value = "40.88"
$ws.Range("E10").Value = "  +10.31%  "
$ws.Range("D11").Value = "0.0929"
$ws.Range("E11").Value = "  +5.95%  "
$ws.Range("D12").Value = "6.89"
$ws.Range("E12").Value = "  +7.16%  "
$ws.Range("E13").Value = "  +2.94%  "
$ws.Range("D14").Value = "2.576.50"
$ws.Range("E14").Value = "  +5.17%  "
$ws.Range("D15").Value = "14.58"
$ws.Range("E15").Value = "  +3.57%  "
$ws.Range("D16").Value = "2.245.01"
$ws.Range("E16").Value = "  +4.40%  "
$ws.Range("D17").Value = "0.791"
$ws.Range("E17").Value = "  +3.79%  "
$ws.Range("D18").Value = "42.953.92"
$ws.Range("E18").Value = "  +5.78%  "
$ws.Range("D19").Value = "0.0000105"
$ws.Range("E19").Value = "  +8.33%  "
$ws.Range("D20").Value = "71.17"
$ws.Range("E20").Value = "  +3.78%  "
$ws.Range("D21").Value = "6.00"
$ws.Range("E21").Value = "  +6.48%  "
$ws.Range("D22").Value = "9.79"
$ws.Range("E22").Value = "  +5.60%  "
$ws.Range("D23").Value = "2.21"
$ws.Range("E23").Value = "  +21.17%  "
$ws.Range("D24").Value = "229.54"
$ws.Range("E24").Value = "  +4.38%  "
$ws.Range("E25").Value = "  +0.27%  "
$ws.Range("E26").Value = "  +4.73%  "
$ws.Range("D27").Value = "3.44"
$ws.Range("E27").Value = "  +7.09%  "
$ws.Range("D28").Value = "39.15"
$ws.Range("E28").Value = "  +33.52%  "
$ws.Range("D29").Value = "2.24"
$ws.Range("E29").Value = "  +6.50%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").Value = "171.37"
$ws.Range("E30").Value = "  +2.84%  "
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").Value = "2.14"
$ws.Range("E31").Value = "  +0.84%  "
$ws.Range("D32").Value = "20.23"
$ws.Range("E32").Value = "  +5.16%  "
$ws.Range("D33").Value = "0.0801"
$ws.Range("E33").Value = "  +9.37%  "
$ws.Range("E34").Value = "  +7.90%  "
$ws.Range("E35").Value = "  +3.89%  "
$ws.Range("D36").Value = "0.109"
$ws.Range("E36").Value = "  +14.50%  "
$ws.Range("D37").Value = "4.46"
$ws.Range("E37").Value = "  +13.25%  "
$ws.Range("E38").Value = "  +22.31%  "
$ws.Range("D39").Value = "12.98"
$ws.Range("E39").Value = "  +17.03%  "
$ws.Range("E40").Value = "  +5.94%  "
$ws.Range("E41").Value = "  +13.95%  "
$ws.Range("D42").Value = "5.40"
$ws.Range("E42").Value = "  +4.31%  "
$ws.Range("D43").Value = "59.20"
$ws.Range("E43").Value = "  +5.97%  "
$ws.Range("D44").Value = "104.42"
$ws.Range("E44").Value = "  +10.83%  "
$ws.Range("D45").Value = "8.68"
$ws.Range("E45").Value = "  +8.06%  "
$ws.Range("D46").Value = "0.480"
$ws.Range("E46").Value = "  +36.63%  "
$ws.Range("D47").Value = "0.0989"
$ws.Range("E47").Value = "  +5.33%  "
$ws.Range("D48").Value = "2.41"
$ws.Range("E48").Value = "  +15.58%  "
$ws.Range("E49").Value = "  +5.90%  "
$ws.Range("E50").Value = "  +6.30%  "
$ws.Range("E51").Value = "  +3.59%  "

# Restore default style on column D (undo the NumberFormat-driven style bump)
$ws.Range("D2:D51").Style = "Normal"
